# Replace all instances of the "Bookerly" font with "Times New Roman"
# so the document renders consistently across platforms that may lack
# the Bookerly font installed.
$d = $word.ActiveDocument

for ($i = 1; $i -le $d.Styles.Count; $i++) {
    $style = $d.Styles.Item($i)
    if ($style.Font.Name -eq "Bookerly") {
        $style.Font.Name = "Times New Roman"
    }
}
